# Update "想去人数" (interested-count) figures after a re-scrape.
# 展览 sheet: 南宁·草莓动漫节 (row 3) 2168 -> 2170
#             南宁·2024三月三国潮动漫节（良牙春典）(row 5) 1466 -> 1480
# 全部类型 sheet carries the same two events in rows 3 and 7.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2170
$wsExhibit.Range("F5").Value = 1480

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2170
$wsAll.Range("F7").Value = 1480
